# Auto-generated script applying 2024-05-04 data update
# Source: commit "Add data for 2024-05-04" diff of output/violent-crime-ytd.xlsx

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2470
$ws.Range("K3").Value = 2376
$ws.Range("K4").Value = 500
$ws.Range("K6").Value = 2963
$ws.Range("K7").Value = 8463

$ws = $wb.Sheets.Item("Logan Square")
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 123

$ws = $wb.Sheets.Item("Austin")
$ws.Range("K2").Value = 167
$ws.Range("K3").Value = 166
$ws.Range("K7").Value = 566

$ws = $wb.Sheets.Item("South Chicago")
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 187

$ws = $wb.Sheets.Item("Garfield Park")
$ws.Range("K2").Value = 93
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 327

$ws = $wb.Sheets.Item("West Pullman")
$ws.Range("K3").Value = 43
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 135

$ws = $wb.Sheets.Item("Grand Crossing")
$ws.Range("K2").Value = 69
$ws.Range("K3").Value = 96
$ws.Range("K7").Value = 271

$ws = $wb.Sheets.Item("New City")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 201

$ws = $wb.Sheets.Item("Woodlawn")
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 152

$ws = $wb.Sheets.Item("By Neighborhood")
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 255
$ws.Range("K8").Value = 566
$ws.Range("K11").Value = 183
$ws.Range("K12").Value = 12
$ws.Range("K14").Value = 48
$ws.Range("K15").Value = 84
$ws.Range("K16").Value = 24
$ws.Range("K18").Value = 56
$ws.Range("K19").Value = 250
$ws.Range("K20").Value = 188
$ws.Range("K23").Value = 78
$ws.Range("K29").Value = 433
$ws.Range("K33").Value = 327
$ws.Range("K37").Value = 271
$ws.Range("K42").Value = 297
$ws.Range("K46").Value = 17
$ws.Range("K47").Value = 47
$ws.Range("K48").Value = 103
$ws.Range("K49").Value = 56
$ws.Range("K50").Value = 54
$ws.Range("K51").Value = 92
$ws.Range("K52").Value = 232
$ws.Range("K53").Value = 123
$ws.Range("K54").Value = 160
$ws.Range("K55").Value = 93
$ws.Range("K63").Value = 32
$ws.Range("K64").Value = 55
$ws.Range("K65").Value = 201
$ws.Range("K74").Value = 10
$ws.Range("K76").Value = 121
$ws.Range("K78").Value = 112
$ws.Range("K79").Value = 220
$ws.Range("K83").Value = 187
$ws.Range("K84").Value = 61
$ws.Range("K85").Value = 410
$ws.Range("K86").Value = 55
$ws.Range("K94").Value = 103
$ws.Range("K95").Value = 135
$ws.Range("K99").Value = 152
$ws.Range("K101").Value = 8463

$ws = $wb.Sheets.Item("South Deering")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 61

$ws = $wb.Sheets.Item("Lincoln Park")
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 56

$ws = $wb.Sheets.Item("Loop")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 160

$ws = $wb.Sheets.Item("Englewood")
$ws.Range("K2").Value = 115
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 433

$ws = $wb.Sheets.Item("Lake View")
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 103

$ws = $wb.Sheets.Item("Chatham")
$ws.Range("K3").Value = 66
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 250

$ws = $wb.Sheets.Item("River North")
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 121

$ws = $wb.Sheets.Item("Bridgeport")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 48

$ws = $wb.Sheets.Item("Ashburn")
$ws.Range("K2").Value = 23
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 67

$ws = $wb.Sheets.Item("Humboldt Park")
$ws.Range("K3").Value = 91
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 297

$ws = $wb.Sheets.Item("Rogers Park")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 112

$ws = $wb.Sheets.Item("Lower West Side")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 93

$ws = $wb.Sheets.Item("Jefferson Park")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 17

$ws = $wb.Sheets.Item("Douglas")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 78

$ws = $wb.Sheets.Item("Roseland")
$ws.Range("K2").Value = 72
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 220

$ws = $wb.Sheets.Item("Near South Side")
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 55

$ws = $wb.Sheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 64
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 188

$ws = $wb.Sheets.Item("Calumet Heights")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 56

$ws = $wb.Sheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 84
$ws.Range("K7").Value = 255

$ws = $wb.Sheets.Item("West Loop")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 103

$ws = $wb.Sheets.Item("Kenwood")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 47

$ws = $wb.Sheets.Item("Brighton Park")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 84

$ws = $wb.Sheets.Item("Lincoln Square")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 54

$ws = $wb.Sheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 183

$ws = $wb.Sheets.Item("Streeterville")
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 55

$ws = $wb.Sheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 25
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 92

$ws = $wb.Sheets.Item("South Shore")
$ws.Range("K2").Value = 150
$ws.Range("K3").Value = 139
$ws.Range("K7").Value = 410

$ws = $wb.Sheets.Item("Little Village")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 232

$ws = $wb.Sheets.Item("Beverly")
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 12

$ws = $wb.Sheets.Item("Bucktown")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 24

$ws = $wb.Sheets.Item("Printers Row")
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 10
